# fb-surveyState.xlsx update: "data up to 12"
# - Fills in previously-missing Puerto Rico (AR) values for rows 100-105
# - Corrects a handful of recomputed percentages in rows 119-122
# - Appends 8 new survey-date rows (123-130) for 01 Jun 2020 .. 08 Jun 2020
#   (row 123 already had its date in column A; this adds the rest of its data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Corrected values in rows 119-122 ---
$modifiedCells = @(
    @("M119", 0.4127641),
    @("M120", 0.3987804),
    @("M121", 0.3850001),
    @("F122", 0.5341384),
    @("G122", 0.3034544),
    @("H122", 0.3358031),
    @("L122", 0.3555895),
    @("M122", 0.3554228),
    @("T122", 0.3411089),
    @("AN122", 0.3508618),
    @("AV122", 0.3254),
    @("AW122", 0.3770781),
    @("BB122", 0.3145865)
)
foreach ($m in $modifiedCells) {
    $ws.Range($m[0]).Value = $m[1]
}

# --- 2) Newly-filled Puerto Rico (AR) values for rows 100-105 ---
$arCells = @(
    @("AR100", 0),
    @("AR101", 0.4651163),
    @("AR102", 0.462963),
    @("AR103", 0.4716981),
    @("AR104", 0.5181346999999999),
    @("AR105", 0.5617978)
)
foreach ($a in $arCells) {
    $ws.Range($a[0]).Value = $a[1]
}

# --- 3) New survey rows 123-130 ---
# Row 123: 01 06 2020
$ws.Range("A123").Value = "01 06 2020"
$row123 = @(
    @("B123", 0.4849849),
    @("C123", 0.4460406),
    @("D123", 0.6138211),
    @("F123", 0.5270899999999999),
    @("G123", 0.3018257),
    @("H123", 0.351763),
    @("I123", 0.3075339),
    @("J123", 0.3026906),
    @("K123", 0.4059549),
    @("L123", 0.3531876),
    @("M123", 0.3403102),
    @("O123", 0.1599468),
    @("P123", 0.5104031),
    @("Q123", 0.4898054),
    @("R123", 0.3781969),
    @("S123", 0.4894694),
    @("T123", 0.3493179),
    @("U123", 0.347692),
    @("V123", 0.4493957),
    @("W123", 0.2881108),
    @("X123", 0.5387116),
    @("Y123", 0.1975282),
    @("Z123", 0.2931544),
    @("AA123", 0.4932034),
    @("AB123", 0.4024588),
    @("AD123", 0.6372597),
    @("AE123", 0.3857645),
    @("AF123", 0.3498461),
    @("AG123", 0.4965635),
    @("AH123", 0.4398901),
    @("AI123", 0.1508344),
    @("AJ123", 0.4169974),
    @("AK123", 0.304897),
    @("AL123", 0.3641499),
    @("AM123", 0.3244266),
    @("AN123", 0.3437967),
    @("AO123", 0.4885348),
    @("AP123", 0.3387037),
    @("AQ123", 0.2826041),
    @("AS123", 0.3010444),
    @("AT123", 0.3556507),
    @("AU123", 0.3893191),
    @("AV123", 0.3718593),
    @("AW123", 0.3853534),
    @("AX123", 0.376092),
    @("AY123", 0.4326245),
    @("BA123", 0.20475),
    @("BB123", 0.3072368),
    @("BC123", 0.3163836),
    @("BD123", 0.4023052),
    @("BE123", 0.5168106)
)
foreach ($c in $row123) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 124: 02 06 2020
$ws.Range("A124").Value = "02 06 2020"
$row124 = @(
    @("B124", 0.4873674),
    @("C124", 0.4277594),
    @("D124", 0.6214668),
    @("F124", 0.5100114),
    @("G124", 0.2929021),
    @("H124", 0.3898445),
    @("I124", 0.3040759),
    @("J124", 0.3033708),
    @("K124", 0.3712185),
    @("L124", 0.340961),
    @("M124", 0.345522),
    @("O124", 0.1708493),
    @("P124", 0.5606184),
    @("Q124", 0.4700609),
    @("R124", 0.3841517),
    @("S124", 0.5279946),
    @("T124", 0.3587215),
    @("U124", 0.3377187),
    @("V124", 0.4609262),
    @("W124", 0.2676403),
    @("X124", 0.4957328),
    @("Y124", 0.242246),
    @("Z124", 0.278411),
    @("AA124", 0.5003625),
    @("AB124", 0.4403837),
    @("AD124", 0.6892334),
    @("AE124", 0.390612),
    @("AF124", 0.3457118),
    @("AG124", 0.5791641),
    @("AH124", 0.5088089),
    @("AI124", 0.1621912),
    @("AJ124", 0.3987549),
    @("AK124", 0.3244725),
    @("AL124", 0.4130027),
    @("AM124", 0.3172603),
    @("AN124", 0.39364),
    @("AO124", 0.5011736),
    @("AP124", 0.2958809),
    @("AQ124", 0.2920578),
    @("AS124", 0.2494945),
    @("AT124", 0.3378541),
    @("AU124", 0.433814),
    @("AV124", 0.3926633),
    @("AW124", 0.4011364),
    @("AX124", 0.4307844),
    @("AY124", 0.4289623),
    @("BA124", 0.2161714),
    @("BB124", 0.2678329),
    @("BC124", 0.3306823),
    @("BD124", 0.4189048),
    @("BE124", 0.493441)
)
foreach ($c in $row124) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 125: 03 06 2020
$ws.Range("A125").Value = "03 06 2020"
$row125 = @(
    @("B125", 0.512557),
    @("C125", 0.3921674),
    @("D125", 0.7087761),
    @("F125", 0.4884358),
    @("G125", 0.2974378),
    @("H125", 0.3925943),
    @("I125", 0.3024579),
    @("J125", 0.3161593),
    @("K125", 0.3373779),
    @("L125", 0.3305049),
    @("M125", 0.35402),
    @("O125", 0.1306482),
    @("P125", 0.5445248),
    @("Q125", 0.4319583),
    @("R125", 0.3705051),
    @("S125", 0.577251),
    @("T125", 0.3550279),
    @("U125", 0.3175135),
    @("V125", 0.4448556),
    @("W125", 0.2625416),
    @("X125", 0.4692786),
    @("Y125", 0.278504),
    @("Z125", 0.2815379),
    @("AA125", 0.4839492),
    @("AB125", 0.4063553),
    @("AD125", 0.6052394),
    @("AE125", 0.365318),
    @("AF125", 0.3746594),
    @("AG125", 0.5343166),
    @("AH125", 0.5091483),
    @("AI125", 0.2034716),
    @("AJ125", 0.3600452),
    @("AK125", 0.3011364),
    @("AL125", 0.4680613),
    @("AM125", 0.3043706),
    @("AN125", 0.3866752),
    @("AO125", 0.4857331),
    @("AP125", 0.2554321),
    @("AQ125", 0.3352552),
    @("AS125", 0.2143884),
    @("AT125", 0.3535209),
    @("AU125", 0.4010737),
    @("AV125", 0.412172),
    @("AW125", 0.4023345),
    @("AX125", 0.4466162),
    @("AY125", 0.423601),
    @("BA125", 0.2345798),
    @("BB125", 0.2754108),
    @("BC125", 0.3200175),
    @("BD125", 0.4065789),
    @("BE125", 0.5604767)
)
foreach ($c in $row125) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 126: 04 06 2020
$ws.Range("A126").Value = "04 06 2020"
$row126 = @(
    @("B126", 0.4131575),
    @("C126", 0.4384974),
    @("D126", 0.7406291),
    @("F126", 0.4806303),
    @("G126", 0.3127366),
    @("H126", 0.3757393),
    @("I126", 0.2860209),
    @("J126", 0.283773),
    @("K126", 0.3585937),
    @("L126", 0.3358304),
    @("M126", 0.3677662),
    @("O126", 0.2038118),
    @("P126", 0.5112910000000001),
    @("Q126", 0.4220976),
    @("R126", 0.3400707),
    @("S126", 0.5947778),
    @("T126", 0.3167692),
    @("U126", 0.3491108),
    @("V126", 0.4859583),
    @("W126", 0.281962),
    @("X126", 0.4860965),
    @("Y126", 0.2322585),
    @("Z126", 0.3001926),
    @("AA126", 0.4645572),
    @("AB126", 0.3897815),
    @("AD126", 0.5196123),
    @("AE126", 0.4404095),
    @("AF126", 0.3783631),
    @("AG126", 0.579039),
    @("AH126", 0.5171445),
    @("AI126", 0.2167148),
    @("AJ126", 0.3161775),
    @("AK126", 0.2759328),
    @("AL126", 0.4497735),
    @("AM126", 0.327866),
    @("AN126", 0.3619751),
    @("AO126", 0.5234121),
    @("AP126", 0.2610234),
    @("AQ126", 0.3351463),
    @("AS126", 0.3588017),
    @("AT126", 0.353308),
    @("AU126", 0.3238842),
    @("AV126", 0.4470904),
    @("AW126", 0.4167661),
    @("AX126", 0.4855886),
    @("AY126", 0.4192885),
    @("BA126", 0.2924271),
    @("BB126", 0.2700101),
    @("BC126", 0.3393596),
    @("BD126", 0.4174135),
    @("BE126", 0.5463456)
)
foreach ($c in $row126) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 127: 05 06 2020
$ws.Range("A127").Value = "05 06 2020"
$row127 = @(
    @("B127", 0.4168694),
    @("C127", 0.4910977),
    @("D127", 0.7540678),
    @("F127", 0.512339),
    @("G127", 0.3046411),
    @("H127", 0.378775),
    @("I127", 0.2867077),
    @("J127", 0.2466189),
    @("K127", 0.3866753),
    @("L127", 0.3347589),
    @("M127", 0.3875201),
    @("O127", 0.1681046),
    @("P127", 0.4705472),
    @("Q127", 0.4456625),
    @("R127", 0.3334852),
    @("S127", 0.5629997),
    @("T127", 0.3910009),
    @("U127", 0.3223967),
    @("V127", 0.5048329),
    @("W127", 0.2622208),
    @("X127", 0.4956116),
    @("Y127", 0.2462375),
    @("Z127", 0.2792438),
    @("AA127", 0.4654325),
    @("AB127", 0.4086622),
    @("AD127", 0.5096232000000001),
    @("AE127", 0.3566518),
    @("AF127", 0.3736788),
    @("AG127", 0.6509494),
    @("AH127", 0.5169411),
    @("AI127", 0.215527),
    @("AJ127", 0.3080879),
    @("AK127", 0.2942504),
    @("AL127", 0.4608879),
    @("AM127", 0.3600704),
    @("AN127", 0.3427985),
    @("AO127", 0.5053021),
    @("AP127", 0.256),
    @("AQ127", 0.3137508),
    @("AS127", 0.3438769),
    @("AT127", 0.3827677),
    @("AU127", 0.2853707),
    @("AV127", 0.4441546),
    @("AW127", 0.3848316),
    @("AX127", 0.4692577),
    @("AY127", 0.4147755),
    @("BA127", 0.2757415),
    @("BB127", 0.2490902),
    @("BC127", 0.2902158),
    @("BD127", 0.425783),
    @("BE127", 0.5246684)
)
foreach ($c in $row127) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 128: 06 06 2020
$ws.Range("A128").Value = "06 06 2020"
$row128 = @(
    @("B128", 0.4614846),
    @("C128", 0.5491541),
    @("D128", 0.7291177),
    @("F128", 0.5220318),
    @("G128", 0.3303167),
    @("H128", 0.4171816),
    @("I128", 0.2901728),
    @("J128", 0.2287582),
    @("K128", 0.418976),
    @("L128", 0.3412584),
    @("M128", 0.3823952),
    @("O128", 0.2062544),
    @("P128", 0.449665),
    @("Q128", 0.3681595),
    @("R128", 0.3655015),
    @("S128", 0.5656745),
    @("T128", 0.4070283),
    @("U128", 0.2958613),
    @("V128", 0.4753994),
    @("W128", 0.2316502),
    @("X128", 0.5120595999999999),
    @("Y128", 0.2399783),
    @("Z128", 0.2983476),
    @("AA128", 0.4789137),
    @("AB128", 0.4574438),
    @("AD128", 0.5611161),
    @("AE128", 0.3015285),
    @("AF128", 0.3613804),
    @("AG128", 0.8000958),
    @("AH128", 0.5292734),
    @("AI128", 0.2204649),
    @("AJ128", 0.3298704),
    @("AK128", 0.3102232),
    @("AL128", 0.5064299),
    @("AM128", 0.3650611),
    @("AN128", 0.3561167),
    @("AO128", 0.4422867),
    @("AP128", 0.2418054),
    @("AQ128", 0.3173486),
    @("AS128", 0.3216968),
    @("AT128", 0.4180432),
    @("AU128", 0.326049),
    @("AV128", 0.4494416),
    @("AW128", 0.3906886),
    @("AX128", 0.5031330000000001),
    @("AY128", 0.4107686),
    @("BA128", 0.3134056),
    @("BB128", 0.2758183),
    @("BC128", 0.2867612),
    @("BD128", 0.3846225),
    @("BE128", 0.6298948)
)
foreach ($c in $row128) {
    $ws.Range($c[0]).Value = $c[1]
}

# Row 129: 07 06 2020
$ws.Range("A129").Value = "07 06 2020"

# Row 130: 08 06 2020
$ws.Range("A130").Value = "08 06 2020"
